$d = $word.ActiveDocument
$find = $d.Content.Find

$pairs = @(
    ,@("0+50=", "18+54=")
    ,@("17+42=", "26-9=")
    ,@("67-42=", "95-60=")
    ,@("45-6=", "13+58=")
    ,@("63+8=", "32+23=")
    ,@("99-56=", "74-13=")
    ,@("7+55=", "86-50=")
    ,@("97-30=", "47+21=")
    ,@("44-2=", "81-16=")
    ,@("53-35=", "26+54=")
    ,@("77-69=", "92-41=")
    ,@("68-65=", "13+71=")
    ,@("75-15=", "76-1=")
    ,@("67-38=", "86-31=")
    ,@("42-27=", "18+53=")
    ,@("69+22=", "35-17=")
    ,@("99-65=", "85+1=")
    ,@("21+57=", "50+26=")
    ,@("38-15=", "59-55=")
    ,@("10+83=", "88-88=")
    ,@("89-45=", "63+31=")
    ,@("45-7=", "92+3=")
    ,@("70-57=", "39+25=")
    ,@("59+6=", "73+4=")
    ,@("26-18=", "95-3=")
    ,@("43-38=", "66+9=")
    ,@("22+32=", "36+19=")
    ,@("39-22=", "81-43=")
    ,@("16+33=", "92-6=")
    ,@("46+32=", "15+5=")
    ,@("73-48=", "82-52=")
    ,@("18+29=", "90-71=")
    ,@("47+49=", "97-91=")
    ,@("57+37=", "34+29=")
    ,@("45-44=", "94-51=")
    ,@("1+66=", "50-10=")
    ,@("97-95=", "8+30=")
    ,@("73-2=", "31-15=")
    ,@("33+27=", "79+15=")
    ,@("51+44=", "77-0=")
    ,@("83-66=", "96-52=")
    ,@("2+14=", "59-58=")
    ,@("42-32=", "66-9=")
    ,@("70+2=", "23+2=")
    ,@("57-23=", "38+1=")
    ,@("70-35=", "85-11=")
    ,@("96-40=", "49+20=")
    ,@("85-59=", "22+68=")
    ,@("96-26=", "34+24=")
    ,@("49-18=", "99-33=")
    ,@("68+0=", "88-39=")
    ,@("85-8=", "88-59=")
    ,@("31+22=", "25+35=")
    ,@("86-84=", "72-13=")
    ,@("65-56=", "3+5=")
    ,@("87-15=", "0+54=")
    ,@("9-5=", "59+14=")
    ,@("77-76=", "80-0=")
    ,@("42-2=", "25+73=")
    ,@("21-18=", "99-74=")
    ,@("4+75=", "22-12=")
    ,@("94-60=", "49+48=")
    ,@("51-33=", "80-78=")
    ,@("64+9=", "27+71=")
    ,@("55-18=", "64+33=")
    ,@("21+75=", "66-13=")
    ,@("27+15=", "3+94=")
    ,@("96-69=", "92-15=")
    ,@("25+5=", "20+39=")
    ,@("32+64=", "52-42=")
    ,@("41+37=", "80+19=")
    ,@("33+46=", "90-49=")
    ,@("72+18=", "14+39=")
    ,@("52-25=", "88-44=")
    ,@("34-15=", "93-39=")
    ,@("17-10=", "53+23=")
    ,@("46-13=", "3+14=")
    ,@("96-34=", "9+35=")
    ,@("10+87=", "43-39=")
    ,@("92-9=", "11+24=")
    ,@("19+24=", "73-32=")
    ,@("10+11=", "75-47=")
    ,@("92-88=", "76+22=")
    ,@("78-50=", "73+14=")
    ,@("29+53=", "90-41=")
    ,@("60-23=", "44+11=")
    ,@("57+6=", "26-2=")
    ,@("39-35=", "7+3=")
    ,@("13+67=", "61-18=")
    ,@("96+3=", "54-48=")
    ,@("3+89=", "17-12=")
    ,@("10+56=", "30-15=")
    ,@("27+0=", "20+41=")
    ,@("68-21=", "5+60=")
    ,@("99-94=", "35-18=")
    ,@("63+17=", "25-25=")
    ,@("74+2=", "13+68=")
    ,@("19+36=", "79+12=")
    ,@("43-12=", "27+53=")
    ,@("54+36=", "32-27=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $result = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "FAILED to replace: $old -> $new"
    }
}

Write-Output "Done"
